$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1030755.25
$ws.Range("I88").Value = 495.5
$ws.Range("J88").Value = 1545885.1
$ws.Range("K88").Value = 495.5
$ws.Range("L88").Value = 1545885.1
$ws.Range("M88").Value = -89.5
$ws.Range("N88").Value = -1546697.1

$ws.Range("H91").Value = 1030755.25
$ws.Range("I91").Value = 495.5
$ws.Range("J91").Value = 1545885.1
$ws.Range("K91").Value = 495.5
$ws.Range("L91").Value = 1545885.1
$ws.Range("M91").Value = 908.5
$ws.Range("N91").Value = -1548693.1

$ws.Range("H125").Value = 1625.4
$ws.Range("I125").Value = 1715.5
$ws.Range("K125").Value = 15439.5
$ws.Range("M125").Value = -12979.5

$ws.Range("H137").Value = 1121.5625
$ws.Range("I137").Value = 1002.2353
$ws.Range("J137").Value = 1256.8
$ws.Range("K137").Value = 3006.7059
$ws.Range("L137").Value = 3770.4
$ws.Range("M137").Value = -456.7058999999999
$ws.Range("N137").Value = -8870.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3102.2307
$ws.Range("I32").Value = 2881.1453
$ws.Range("K32").Value = 2881.1453
$ws.Range("M32").Value = -2594.1453

$ws.Range("H74").Value = 904.6905
$ws.Range("I74").Value = 502.51724
$ws.Range("J74").Value = 1801.8462
$ws.Range("K74").Value = 502.51724
$ws.Range("L74").Value = 1801.8462
$ws.Range("M74").Value = 371.48276
$ws.Range("N74").Value = -3549.8462

$ws.Range("H77").Value = 904.6905
$ws.Range("I77").Value = 502.51724
$ws.Range("J77").Value = 1801.8462
$ws.Range("K77").Value = 2512.5862
$ws.Range("L77").Value = 9009.231
$ws.Range("M77").Value = 1855.4138
$ws.Range("N77").Value = -17745.231

$ws.Range("H88").Value = 2762.3125
$ws.Range("I88").Value = 2374
$ws.Range("J88").Value = 2817.7856
$ws.Range("K88").Value = 2374
$ws.Range("L88").Value = 2817.7856
$ws.Range("M88").Value = -1968
$ws.Range("N88").Value = -3629.7856

$ws.Range("H91").Value = 2762.3125
$ws.Range("I91").Value = 2374
$ws.Range("J91").Value = 2817.7856
$ws.Range("K91").Value = 2374
$ws.Range("L91").Value = 2817.7856
$ws.Range("M91").Value = -970
$ws.Range("N91").Value = -5625.7856

$ws.Range("H102").Value = 11113452
$ws.Range("I102").Value = 11906984
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 11906984
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -11905362
$ws.Range("N102").Value = -7244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9118.9375
$ws.Range("I134").Value = 1336.7273
$ws.Range("K134").Value = 4010.1819
$ws.Range("M134").Value = -1475.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1523.7812
$ws.Range("I31").Value = 1127.579
$ws.Range("K31").Value = 1127.579
$ws.Range("M31").Value = -832.579

$ws.Range("H34").Value = 1523.7812
$ws.Range("I34").Value = 1127.579
$ws.Range("K34").Value = 1127.579
$ws.Range("M34").Value = -925.579

$ws.Range("H62").Value = 16669371
$ws.Range("I62").Value = 3056.25
$ws.Range("J62").Value = 50002000
$ws.Range("K62").Value = 3056.25
$ws.Range("L62").Value = 50002000
$ws.Range("M62").Value = -2432.25
$ws.Range("N62").Value = -50003248

$ws.Range("H65").Value = 16669371
$ws.Range("I65").Value = 3056.25
$ws.Range("J65").Value = 50002000
$ws.Range("K65").Value = 15281.25
$ws.Range("L65").Value = 250010000
$ws.Range("M65").Value = -12161.25
$ws.Range("N65").Value = -250016240

$ws.Range("H132").Value = 6457.32
$ws.Range("I132").Value = 8152.3125
$ws.Range("J132").Value = 3444
$ws.Range("K132").Value = 24456.9375
$ws.Range("L132").Value = 10332
$ws.Range("M132").Value = -21926.9375
$ws.Range("N132").Value = -15392

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6251729.5
$ws.Range("I34").Value = 1620
$ws.Range("J34").Value = 8335099
$ws.Range("K34").Value = 4860
$ws.Range("L34").Value = 25005297
$ws.Range("M34").Value = -4776
$ws.Range("N34").Value = -25005465

$ws.Range("H39").Value = 1916.64
$ws.Range("J39").Value = 1659.8182
$ws.Range("L39").Value = 4979.4546
$ws.Range("N39").Value = -5567.4546

$ws.Range("H55").Value = 2500.8333
$ws.Range("J55").Value = 3251.25
$ws.Range("L55").Value = 9753.75
$ws.Range("N55").Value = -10107.75

$ws.Range("H63").Value = 5232
$ws.Range("I63").Value = 2941.3333
$ws.Range("K63").Value = 8823.999899999999
$ws.Range("M63").Value = -8074.999899999999

$ws.Range("H66").Value = 5232
$ws.Range("I66").Value = 2941.3333
$ws.Range("K66").Value = 26471.9997
$ws.Range("M66").Value = -22727.9997

$ws.Range("H70").Value = 3455
$ws.Range("I70").Value = 1207
$ws.Range("J70").Value = 4391.6665
$ws.Range("K70").Value = 3621
$ws.Range("L70").Value = 13174.9995
$ws.Range("M70").Value = -3306
$ws.Range("N70").Value = -13804.9995

$ws.Range("H73").Value = 3455
$ws.Range("I73").Value = 1207
$ws.Range("J73").Value = 4391.6665
$ws.Range("K73").Value = 3621
$ws.Range("L73").Value = 13174.9995
$ws.Range("M73").Value = -2529
$ws.Range("N73").Value = -15358.9995

$ws.Range("H92").Value = 789.9
$ws.Range("I92").Value = 839.2
$ws.Range("J92").Value = 740.6
$ws.Range("K92").Value = 2517.6
$ws.Range("L92").Value = 2221.8
$ws.Range("M92").Value = -1269.6
$ws.Range("N92").Value = -4717.8

$ws.Range("H93").Value = 5261.4546
$ws.Range("J93").Value = 5261.4546
$ws.Range("L93").Value = 15784.3638
$ws.Range("N93").Value = -19528.3638

$ws.Range("H129").Value = 32052050
$ws.Range("I129").Value = 37037500
$ws.Range("J129").Value = 20834792
$ws.Range("K129").Value = 111112500
$ws.Range("L129").Value = 62504376
$ws.Range("M129").Value = -111107500
$ws.Range("N129").Value = -62514376

$ws.Range("H131").Value = 11237121
$ws.Range("J131").Value = 1219.262
$ws.Range("L131").Value = 3657.786
$ws.Range("N131").Value = -13737.786

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2172.182
$ws.Range("I126").Value = 1760
$ws.Range("K126").Value = 5280
$ws.Range("M126").Value = -2810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 240
$ws.Range("I55").Value = 176.85715
$ws.Range("J55").Value = 271.57144
$ws.Range("K55").Value = 176.85715
$ws.Range("L55").Value = 271.57144
$ws.Range("M55").Value = -3.85714999999999
$ws.Range("N55").Value = -617.5714399999999

$ws.Range("H93").Value = 1176.7778
$ws.Range("I93").Value = 665.25
$ws.Range("K93").Value = 665.25
$ws.Range("M93").Value = 582.75

$ws.Range("H122").Value = 8069197.5
$ws.Range("I122").Value = 10875269
$ws.Range("K122").Value = 32625807
$ws.Range("M122").Value = -32623357

$ws.Range("H128").Value = 99990
$ws.Range("J128").Value = 99990
$ws.Range("L128").Value = 99990
$ws.Range("N128").Value = -109950

$ws.Range("H136").Value = 5439.591
$ws.Range("I136").Value = 8014.2856
$ws.Range("J136").Value = 933.875
$ws.Range("K136").Value = 24042.8568
$ws.Range("L136").Value = 2801.625
$ws.Range("M136").Value = -21492.8568
$ws.Range("N136").Value = -7901.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 58824492
$ws.Range("I126").Value = 71429420
$ws.Range("J126").Value = 1466.6666
$ws.Range("K126").Value = 214288260
$ws.Range("L126").Value = 4399.9998
$ws.Range("M126").Value = -214285790
$ws.Range("N126").Value = -9339.9998
